# "Updated commit on may 02" - refresh the may-02 exam/schedule/location
# names on the active STAGE sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters here only insofar as it controls the order new entries
# land in the shared-string table; update Schedule, then Location, then
# Exam name to match the source ordering.
$ws.Range("E2").Value = "Pipeline Schedule may 02"
$ws.Range("L2").Value = "LOC mar 0502"
$ws.Range("A2").Value = "SmokeTest Exam may 02"

# Leave the cursor parked on B10, as the author did before committing.
[void]$ws.Range("B10").Select()
